$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.386.35'
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.811.16'
$ws.Range("E3").Value = '  -0.82%  '
$ws.Range("E4").Value = '  -0.34%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.77'
$ws.Range("E5").Value = '  -1.14%  '
$ws.Range("E6").Value = '  -0.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5136'
$ws.Range("E7").Value = '  -0.49%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3999'
$ws.Range("E8").Value = '  +3.51%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07867'
$ws.Range("E9").Value = '  -5.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.115'
$ws.Range("E10").Value = '  -0.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '40.99'
$ws.Range("E11").Value = '  -2.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.380'
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.41'
$ws.Range("E14").Value = '  -3.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.330'
$ws.Range("E15").Value = '  -2.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.808.20'
$ws.Range("E16").Value = '  -1.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.80'
$ws.Range("E17").Value = '  -1.26%  '
$ws.Range("E18").Value = '  -3.54%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06577'
$ws.Range("E19").Value = '  -0.91%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  -0.30%  '
$ws.Range("E21").Value = '  -2.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.019'
$ws.Range("E22").Value = '  -0.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.437.13'
$ws.Range("E23").Value = '  -0.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.15'
$ws.Range("E24").Value = '  -2.59%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.233'
$ws.Range("E25").Value = '  -0.51%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.84'
$ws.Range("E26").Value = '  +1.56%  '
$ws.Range("E27").Value = '  -2.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.019.78'
$ws.Range("E28").Value = '  -0.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.412'
$ws.Range("E29").Value = '  -0.25%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.28'
$ws.Range("E30").Value = '  +1.94%  '
$ws.Range("E32").Value = '  -2.61%  '
$ws.Range("E33").Value = '  -0.45%  '
$ws.Range("E34").Value = '  -2.59%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07255'
$ws.Range("E35").Value = '  -5.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.236'
$ws.Range("E36").Value = '  +5.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02344'
$ws.Range("E37").Value = '  -1.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2175'
$ws.Range("E38").Value = '  -2.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.68'
$ws.Range("E39").Value = '  -2.82%  '
$ws.Range("E40").Value = '  -3.88%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6201'
$ws.Range("E41").Value = '  -3.50%  '
$ws.Range("E42").Value = '  -0.28%  '
$ws.Range("E43").Value = '  -2.85%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.24'
$ws.Range("E44").Value = '  -3.12%  '
$ws.Range("B45").Value = 'WEMIXTOKEN'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.316'
$ws.Range("E45").Value = '  -6.05%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5992'
$ws.Range("E46").Value = '  -3.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.737'
$ws.Range("E47").Value = '  -1.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.84'
$ws.Range("E48").Value = '  -1.62%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.223'
$ws.Range("E49").Value = '  +1.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.932'
$ws.Range("E50").Value = '  -3.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06850'
$ws.Range("E51").Value = '  -1.81%  '
